$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Stage the existing "2011-12-20" string (currently sitting in column I, the
# property_category-adjacent date column) onto the clipboard *before* any cell is
# touched, so it can be replayed verbatim into the new "date" column (J) for every
# row without Excel re-parsing the text and turning it into a date serial number.
$ws.Range("I2").Copy() | Out-Null

for ($r = 2; $r -le 11; $r++) {
    # Capture the values that are about to move one column to the right, before they
    # get overwritten.
    $legislatorId = $ws.Cells.Item($r, 11).Value2
    $idx = $ws.Cells.Item($r, 1).Value2

    # J (col 10) = date ("2011-12-20"), pasted as a value from the clipboard so the
    # text isn't reinterpreted as a date literal.
    $ws.Cells.Item($r, 10).PasteSpecial(-4163) | Out-Null

    # I (col 9) = new "category" column.
    $ws.Cells.Item($r, 9).Value = "normal"

    # K (col 11) = legislator_name (was previously in column J).
    $ws.Cells.Item($r, 11).Value = "黃偉哲"

    # L (col 12) = legislator_id (was previously in column K, numeric).
    $ws.Cells.Item($r, 12).Value = $legislatorId

    # M (col 13) = new "source_file" column.
    $ws.Cells.Item($r, 13).Value = "tmp85f1"

    # N (col 14) = new "index" column, mirrors column A's numeric id.
    $ws.Cells.Item($r, 14).Value = $idx
}

# Header row (row 1): insert "category" header, and append the three trailing
# headers (legislator_id, source_file, index), copying the existing header style
# (bold + border) from K1 onto each new header cell.
$ws.Cells.Item(1, 9).Value = "category"

$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").Value = "legislator_id"
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").Value = "source_file"
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$ws.Range("N1").Value = "index"
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
